$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'307.82"
$ws.Range("E2").Value = "'-4.22%"
$ws.Range("D3").Value = "'39.52"
$ws.Range("E3").Value = "'-7.92%"
$ws.Range("D4").Value = "'5.099"
$ws.Range("E4").Value = "'-1.90%"
$ws.Range("D5").Value = "'0.07690"
$ws.Range("E5").Value = "'-5.93%"
$ws.Range("E6").Value = "'-1.94%"
$ws.Range("D7").Value = "'1.604"
$ws.Range("E7").Value = "'-11.50%"
$ws.Range("D8").Value = "'0.8989"
$ws.Range("E8").Value = "'-3.81%"
$ws.Range("D9").Value = "'0.1006"
$ws.Range("E9").Value = "'-9.11%"
$ws.Range("D10").Value = "'0.1736"
$ws.Range("E10").Value = "'-6.44%"
$ws.Range("D11").Value = "'0.08989"
$ws.Range("E11").Value = "'-5.24%"
$ws.Range("D12").Value = "'0.04418"
$ws.Range("E12").Value = "'-4.73%"
$ws.Range("E13").Value = "'-0.39%"
$ws.Range("D14").Value = "'0.001263"
$ws.Range("E14").Value = "'-2.22%"
$ws.Range("B15").Value = "'TigerCash"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.005890"
$ws.Range("E15").Value = "'0.54%"
$ws.Range("B16").Value = "'UpBots"
$ws.Range("C16").Value = "'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$ws.Range("D16").Value = "'0.007491"
$ws.Range("E16").Value = "'2,413.02%"
$ws.Range("B17").Value = "'LEO"
$ws.Range("C17").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.353"
$ws.Range("E17").Value = "'-0.05%"
$ws.Range("B18").Value = "'BTSEToken"
$ws.Range("C18").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").Value = "'2.419"
$ws.Range("E18").Value = "'-3.53%"
$ws.Range("B19").Value = "'BitpandaEcosystemToken"
$ws.Range("C19").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").Value = "'0.3315"
$ws.Range("E19").Value = "'-1.00%"
$ws.Range("B20").Value = "'MCDex"
$ws.Range("C20").Value = "'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D20").Value = "'7.068"
$ws.Range("E20").Value = "'-5.47%"
$ws.Range("B21").Value = "'ProBitToken"
$ws.Range("C21").Value = "'https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1349"
$ws.Range("E21").Value = "'-2.27%"
$ws.Range("B22").Value = "'ZBToken"
$ws.Range("C22").Value = "'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").Value = "'0.2763"
$ws.Range("E22").Value = "'9.55%"
$ws.Range("B23").Value = "'CoinExToken"
$ws.Range("C23").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").Value = "'0.04147"
$ws.Range("E23").Value = "'-0.03%"
$ws.Range("D24").Value = "'0.001210"
$ws.Range("E24").Value = "'-2.65%"
$ws.Range("D25").Value = "'0.004058"
$ws.Range("E25").Value = "'-5.48%"
$ws.Range("D26").Value = "'0.0001304"
$ws.Range("E26").Value = "'8.59%"
$ws.Range("D38").Value = "'0.02337"
$ws.Range("E38").Value = "'-13.84%"
$ws.Range("D39").Value = "'0.05152"
$ws.Range("E39").Value = "'-7.18%"
$ws.Range("D40").Value = "'0.007935"
$ws.Range("E40").Value = "'-1.83%"
$ws.Range("D41").Value = "'0.1320"
$ws.Range("E41").Value = "'-5.54%"
$ws.Range("D42").Value = "'0.006303"
$ws.Range("E42").Value = "'-3.75%"
$ws.Range("D43").Value = "'0.001956"
$ws.Range("E43").Value = "'-4.21%"
$ws.Range("D44").Value = "'0.008237"
$ws.Range("E44").Value = "'-0.20%"
$ws.Range("D45").Value = "'0.3335"
$ws.Range("E45").Value = "'-4.70%"
$ws.Range("D46").Value = "'0.00006523"
$ws.Range("E46").Value = "'-5.82%"
$ws.Range("E47").Value = "'0.24%"
$ws.Range("E48").Value = "'98.34%"
$ws.Range("D49").Value = "'0.003466"
$ws.Range("E49").Value = "'3.81%"
$ws.Range("D50").Value = "'0.00002107"
$ws.Range("E50").Value = "'0.24%"
$ws.Range("D51").Value = "'0.0002006"
$ws.Range("E51").Value = "'0.24%"
